# Apply the "completed account-service impl and mart-cart-service impl" edit.
#
# Summary of the change (per the OOXML diff):
#  - account-service (row 12) and cart-service (row 14) "Implement" column
#    (G) flips from "Working" to "OK".
#  - A new column I ("Expected to be completed on 5/17/2020") is added next
#    to every micro-service row that still has an outstanding item
#    (rows 12, 14-19), widened to fit the text.
#  - The sheet's selection moves to F25.
#  - Page setup is switched to portrait orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# account-service and cart-service: "Implement" is now done.
$ws.Range("G12").Value = "OK"
$ws.Range("G14").Value = "OK"

# New note column (I) for the services that are still pending completion.
$note = "Expected to be completed on 5/17/2020"
$ws.Range("I12").Value = $note
$ws.Range("I14").Value = $note
$ws.Range("I15").Value = $note
$ws.Range("I16").Value = $note
$ws.Range("I17").Value = $note
$ws.Range("I18").Value = $note
$ws.Range("I19").Value = $note

# Widen the new column so the note is readable (target stored width 34.875
# "characters"; the host quantizes ColumnWidth writes to whole pixels under
# the hood, so 34.14 is the input that lands closest to that value).
$ws.Columns(9).ColumnWidth = 34.14

# Page is printed portrait.
$ws.PageSetup.Orientation = 1

# Leave the selection where the author left it.
$ws.Range("F25").Select() | Out-Null
